# Break out stock.yaml completed
# 1) Normalize the "bsecode" (column D) values on every timeframe sheet from
#    text (inline string) to genuine numbers.
# 2) Append four freshly scraped rows (MRF, BOSCHLTD, MARUTI, GMRINFRA) to the
#    "day" sheet, keeping their bsecode values as text (as produced by the
#    scraper before the later normalization pass).

$wb = $excel.ActiveWorkbook

$sheetNames = @("day", "week", "month", "quarter")
foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $lastRow = $ws.UsedRange.Rows.Count
    for ($r = 2; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, 4)
        $txt = $cell.Value2
        $cell.Value = [double]$txt
    }
}

$daySheet = $wb.Worksheets.Item("day")

$newRows = @(
    @{A=1; B="MRF"; C="Mrf Limited"; D="500290"; E=-0.25; F=126963.3; G=14250; H="day"; I="10/06/2024 10:32:46"},
    @{A=2; B="BOSCHLTD"; C="Bosch Limited"; D="500530"; E=0.17; F=30650; G=12586; H="day"; I="10/06/2024 10:32:46"},
    @{A=3; B="MARUTI"; C="Maruti Suzuki India Limited"; D="532500"; E=-0.73; F=12717.55; G=293471; H="day"; I="10/06/2024 10:32:46"},
    @{A=4; B="GMRINFRA"; C="Gmr Infrastructure Limited"; D="532754"; E=0.27; F=86.93000000000001; G=26759792; H="day"; I="10/06/2024 10:32:46"}
)

$r = $daySheet.UsedRange.Rows.Count + 1
foreach ($row in $newRows) {
    $daySheet.Cells.Item($r, 1).Value = $row.A
    $daySheet.Cells.Item($r, 2).Value = $row.B
    $daySheet.Cells.Item($r, 3).Value = $row.C

    $dcell = $daySheet.Cells.Item($r, 4)
    $dcell.NumberFormat = "@"
    $dcell.Value = $row.D

    $daySheet.Cells.Item($r, 5).Value = $row.E
    $daySheet.Cells.Item($r, 6).Value = $row.F
    $daySheet.Cells.Item($r, 7).Value = $row.G
    $daySheet.Cells.Item($r, 8).Value = $row.H
    $daySheet.Cells.Item($r, 9).Value = $row.I

    $r = $r + 1
}

Write-Host "Updated bsecode columns to numeric and appended 4 rows to 'day' sheet."
